# The ProcessedData sheet's three columns (A, B, C) need to be rotated one
# position to the right: new A = old C, new B = old A, new C = old B.
# This affects the header row (A1:C1) as well as the two data rows
# (A2:C3). Use Copy/PasteSpecial (rather than reading/writing .Value,
# which re-parses numeric-looking text into numbers) so the original
# shared-string cell values and types are preserved exactly, just moved
# to their new columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Stash column C out of the way in a scratch column first, then shift
# B -> C, A -> B, and finally the stashed old C -> A.
$ws.Range("C1:C3").Copy()
$ws.Range("Z1:Z3").PasteSpecial()

$ws.Range("B1:B3").Copy()
$ws.Range("C1:C3").PasteSpecial()

$ws.Range("A1:A3").Copy()
$ws.Range("B1:B3").PasteSpecial()

$ws.Range("Z1:Z3").Copy()
$ws.Range("A1:A3").PasteSpecial()

# Clean up the scratch column so it doesn't linger in the saved sheet.
$ws.Range("Z1:Z3").ClearContents()
